# Move test to tests directory and clean up
# Remove the two extra questionnaire rows (contact details / year founded)
# that were left over from the old copy of this test file, and restore the
# viewport's saved selection so the file doesn't carry stray scroll state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Company")

# Rows 3 and 4 ("Provide the contact details..." / "Please specify the year
# your company was founded") are no longer part of this sample - clear them
# out entirely (values + the now-unused shared strings they referenced).
$ws.Range("A3:F3").ClearContents()
$ws.Range("A4:F4").ClearContents()

# Row 3 had an explicit taller height to fit its (now removed) long question
# text; auto-fit it back down so it reverts to the sheet's default height.
$ws.Rows.Item(3).AutoFit()

# Restore the saved cursor position on the sheet.
[void]$ws.Range("E14").Select()
